# Update Leve price/profit figures across the FFXIV crafting-job profit sheets.
# Source data: scheduled market-board pull (Valefor data centre).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2900
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 3166.6667
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 3166.6667
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -4414.6667
# Row 65
$ws.Range("H65").Value = 2900
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 3166.6667
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 15833.3335
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -22073.3335
# Row 98
$ws.Range("H98").Value = 1842.64
$ws.Range("I98").Value = 1426.3334
$ws.Range("J98").Value = 2467.1
$ws.Range("K98").Value = 1426.3334
$ws.Range("L98").Value = 2467.1
$ws.Range("M98").Value = 71.66660000000002
$ws.Range("N98").Value = -5463.1
# Row 122
$ws.Range("H122").Value = 1842.64
$ws.Range("I122").Value = 1426.3334
$ws.Range("J122").Value = 2467.1
$ws.Range("K122").Value = 4279.0002
$ws.Range("L122").Value = 7401.299999999999
$ws.Range("M122").Value = -1829.0002
$ws.Range("N122").Value = -12301.3
# Row 132
$ws.Range("H132").Value = 2734735
$ws.Range("I132").Value = 4186062.5
$ws.Range("J132").Value = 2824.4119
$ws.Range("K132").Value = 12558187.5
$ws.Range("L132").Value = 8473.235700000001
$ws.Range("M132").Value = -12555657.5
$ws.Range("N132").Value = -13533.2357

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 22459
$ws.Range("J43").Value = 22459
$ws.Range("L43").Value = 22459
$ws.Range("N43").Value = -23085
# Row 109
$ws.Range("H109").Value = 17688.5
$ws.Range("J109").Value = 17688.5
$ws.Range("L109").Value = 17688.5
$ws.Range("N109").Value = -20462.5
# Row 132
$ws.Range("H132").Value = 2580.027
$ws.Range("I132").Value = 1458.15
$ws.Range("J132").Value = 3899.8823
$ws.Range("K132").Value = 4374.450000000001
$ws.Range("L132").Value = 11699.6469
$ws.Range("M132").Value = -1844.450000000001
$ws.Range("N132").Value = -16759.6469

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9806605
$ws.Range("I31").Value = 20834238
$ws.Range("J31").Value = 4264.5186
$ws.Range("K31").Value = 20834238
$ws.Range("L31").Value = 4264.5186
$ws.Range("M31").Value = -20833943
$ws.Range("N31").Value = -4854.5186
# Row 34
$ws.Range("H34").Value = 9806605
$ws.Range("I34").Value = 20834238
$ws.Range("J34").Value = 4264.5186
$ws.Range("K34").Value = 20834238
$ws.Range("L34").Value = 4264.5186
$ws.Range("M34").Value = -20834036
$ws.Range("N34").Value = -4668.5186
# Row 80
$ws.Range("H80").Value = 17000
$ws.Range("J80").Value = 17000
$ws.Range("L80").Value = 17000
$ws.Range("N80").Value = -19246
# Row 83
$ws.Range("H83").Value = 17000
$ws.Range("J83").Value = 17000
$ws.Range("L83").Value = 51000
$ws.Range("N83").Value = -62232
# Row 94
$ws.Range("H94").Value = 2611.2354
$ws.Range("I94").Value = 4304
$ws.Range("J94").Value = 2248.5
$ws.Range("K94").Value = 4304
$ws.Range("L94").Value = 2248.5
$ws.Range("M94").Value = -3853
$ws.Range("N94").Value = -3150.5
# Row 97
$ws.Range("H97").Value = 31997.5
$ws.Range("J97").Value = 31997.5
$ws.Range("L97").Value = 31997.5
$ws.Range("N97").Value = -33979.5
# Row 107
$ws.Range("H107").Value = 1277.8182
$ws.Range("I107").Value = 550.8571
$ws.Range("J107").Value = 2550
$ws.Range("K107").Value = 550.8571
$ws.Range("L107").Value = 2550
$ws.Range("M107").Value = 1369.1429
$ws.Range("N107").Value = -6390
# Row 122
$ws.Range("H122").Value = 27253
$ws.Range("I122").Value = 100012
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 300036
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -297586
$ws.Range("N122").Value = -13900

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 18865.334
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 18865.334
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 18865.334
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -19441.334
# Row 80
$ws.Range("H80").Value = 102501
$ws.Range("I80").Value = 2301.25
$ws.Range("J80").Value = 169300.83
$ws.Range("K80").Value = 2301.25
$ws.Range("L80").Value = 169300.83
$ws.Range("M80").Value = -1303.25
$ws.Range("N80").Value = -171296.83
# Row 81
$ws.Range("H81").Value = 18865.334
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 18865.334
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 18865.334
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -20861.334
# Row 83
$ws.Range("H83").Value = 102501
$ws.Range("I83").Value = 2301.25
$ws.Range("J83").Value = 169300.83
$ws.Range("K83").Value = 11506.25
$ws.Range("L83").Value = 846504.1499999999
$ws.Range("M83").Value = -6514.25
$ws.Range("N83").Value = -856488.1499999999
# Row 84
$ws.Range("H84").Value = 18865.334
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 18865.334
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 56596.00199999999
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -66580.002
# Row 132
$ws.Range("H132").Value = 6516.625
$ws.Range("I132").Value = 9746.667
$ws.Range("J132").Value = 3286.5833
$ws.Range("K132").Value = 29240.001
$ws.Range("L132").Value = 9859.749899999999
$ws.Range("M132").Value = -26710.001
$ws.Range("N132").Value = -14919.7499
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1850.862
$ws.Range("I61").Value = 1174.8667
$ws.Range("J61").Value = 2575.1428
$ws.Range("K61").Value = 1174.8667
$ws.Range("L61").Value = 2575.1428
$ws.Range("M61").Value = -972.8667
$ws.Range("N61").Value = -2979.1428
# Row 80
$ws.Range("H80").Value = 34778
$ws.Range("J80").Value = 34778
$ws.Range("L80").Value = 34778
$ws.Range("N80").Value = -37024
# Row 83
$ws.Range("H83").Value = 34778
$ws.Range("J83").Value = 34778
$ws.Range("L83").Value = 104334
$ws.Range("N83").Value = -115566
# Row 93
$ws.Range("H93").Value = 2573.1052
$ws.Range("I93").Value = 2192.0715
$ws.Range("J93").Value = 3640
$ws.Range("K93").Value = 2192.0715
$ws.Range("L93").Value = 3640
$ws.Range("M93").Value = -944.0715
$ws.Range("N93").Value = -6136
# Row 113
$ws.Range("H113").Value = 1850.862
$ws.Range("I113").Value = 1174.8667
$ws.Range("J113").Value = 2575.1428
$ws.Range("K113").Value = 1174.8667
$ws.Range("L113").Value = 2575.1428
$ws.Range("M113").Value = 995.1333
$ws.Range("N113").Value = -6915.1428
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 413
$ws.Range("I107").Value = 434.6875
$ws.Range("J107").Value = 381.45456
$ws.Range("K107").Value = 1304.0625
$ws.Range("L107").Value = 1144.36368
$ws.Range("M107").Value = 615.9375
$ws.Range("N107").Value = -4984.36368
